$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 53
$ws.Range("C53").Value = 141681
$ws.Range("E53").Value = 590068269

# Row 83
$ws.Range("C83").Value = 3414
$ws.Range("E83").Value = 115758395

# Row 91
$ws.Range("C91").Value = 151142
$ws.Range("E91").Value = 482442285

# Row 92
$ws.Range("C92").Value = 409155
$ws.Range("E92").Value = 1595512102

# Row 95
$ws.Range("C95").Value = 50776
$ws.Range("E95").Value = 932943241

# Row 96
$ws.Range("C96").Value = 17297
$ws.Range("E96").Value = 794564257

# Row 104
$ws.Range("C104").Value = 135245
$ws.Range("E104").Value = 272235386

# Row 116
$ws.Range("C116").Value = 4563
$ws.Range("E116").Value = 20627806

# Row 174
$ws.Range("C174").Value = 226098
$ws.Range("E174").Value = 900654909

# Row 177
$ws.Range("C177").Value = 14719
$ws.Range("E177").Value = 251626957
